# Generate Report for Handoff
# Updates the Priority ("low" -> "ht") and Latest Handoff Datetime values
# for the rows that were "Ready for handoff" (rows 4-7) on both the
# zh-cn and de-de localization-status worksheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4:H7").Value = "2016-08-29 14:35:40"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4:H7").Value = "2016-08-29 14:35:45"

# The Overview sheet's "Latest HO Xliff Generate Date" column mirrors the
# de-de handoff datetime for these rows, so it needs the same update.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G4:G7").Value = "2016-08-29 14:35:45"
